$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.914.44'
$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("D3").Value = '2.615.88'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'309.19"
$ws.Range("E5").Value = '  -1.60%  '
$ws.Range("D6").Value = "'98.28"
$ws.Range("E6").Value = '  -4.18%  '
$ws.Range("E7").Value = '  -1.31%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = "'0.578"
$ws.Range("E9").Value = '  -1.98%  '
$ws.Range("D10").Value = "'38.60"
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").Value = "'54.12"
$ws.Range("E11").Value = '  -1.02%  '
$ws.Range("D12").Value = "'0.0841"
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("E13").Value = '  -3.82%  '
$ws.Range("D14").Value = '3.016.79'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("D16").Value = '2.620.03'
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").Value = "'0.914"
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").Value = "'14.80"
$ws.Range("E18").Value = '  -2.22%  '
$ws.Range("D19").Value = '45.931.92'
$ws.Range("E19").Value = '  -1.91%  '
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("D22").Value = "'12.70"
$ws.Range("E22").Value = '  -5.07%  '
$ws.Range("D23").Value = "'74.50"
$ws.Range("E23").Value = '  +4.66%  '
$ws.Range("D24").Value = "'281.61"
$ws.Range("E24").Value = '  +10.02%  '
$ws.Range("E25").Value = '  -2.51%  '
$ws.Range("E26").Value = '  +1.55%  '
$ws.Range("D27").Value = "'29.48"
$ws.Range("E27").Value = '  +3.28%  '
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").Value = "'4.05"
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("D31").Value = "'38.60"
$ws.Range("E31").Value = '  -7.93%  '
$ws.Range("E32").Value = '  -3.94%  '
$ws.Range("D33").Value = "'6.21"
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = "'3.60"
$ws.Range("E34").Value = '  -4.10%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'2.29"
$ws.Range("E35").Value = '  +0.39%  '
$ws.Range("D36").Value = "'156.63"
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("D38").Value = "'2.81"
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("D39").Value = "'0.122"
$ws.Range("E39").Value = '  +2.75%  '
$ws.Range("D40").Value = "'0.123"
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("E41").Value = '  -7.53%  '
$ws.Range("D42").Value = "'22.09"
$ws.Range("E42").Value = '  +3.41%  '
$ws.Range("D43").Value = "'0.0326"
$ws.Range("E43").Value = '  -1.15%  '
$ws.Range("E44").Value = '  -3.51%  '
$ws.Range("E45").Value = '  -6.95%  '
$ws.Range("D46").Value = '2.104.29'
$ws.Range("E46").Value = '  +3.29%  '
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").Value = "'93.57"
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("D49").Value = "'109.91"
$ws.Range("E49").Value = '  -3.29%  '
$ws.Range("D50").Value = "'9.14"
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("D51").Value = '2.871.19'
$ws.Range("E51").Value = '  -0.21%  '
